# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap displayed country names (shared-string reorder in the source diff) ---
# Rows 197/198: Macao <-> Curazao
$ws.Cells.Item(197, 1).Value = "Curazao"
$ws.Cells.Item(198, 1).Value = "Macao"

# Rows 214/215: Islas Malvinas <-> Montserrat
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 06:09"

# --- Update statistic counters ---
# Row 6 (India)
$ws.Cells.Item(6, 2).Value = 3167323
$ws.Cells.Item(6, 3).Value = 2442
$ws.Cells.Item(6, 4).Value = 2404585
$ws.Cells.Item(6, 5).Value = 704192

# Row 72 (Australia)
$ws.Cells.Item(72, 2).Value = 25053
$ws.Cells.Item(72, 3).Value = 137
$ws.Cells.Item(72, 5).Value = 4925

# Row 122 (Tailandia)
$ws.Cells.Item(122, 2).Value = 3402
$ws.Cells.Item(122, 3).Value = 5
$ws.Cells.Item(122, 4).Value = 3229
$ws.Cells.Item(122, 5).Value = 115

# Row 190 (Butan)
$ws.Cells.Item(190, 2).Value = 156
$ws.Cells.Item(190, 3).Value = 1
$ws.Cells.Item(190, 4).Value = 117
$ws.Cells.Item(190, 5).Value = 39

# Row 197 (now Curazao)
$ws.Cells.Item(197, 2).Value = 47
$ws.Cells.Item(197, 3).Value = 4
$ws.Cells.Item(197, 4).Value = 34
$ws.Cells.Item(197, 5).Value = 12
$ws.Cells.Item(197, 8).Value = 1

# Row 198 (now Macao)
$ws.Cells.Item(198, 2).Value = 46
$ws.Cells.Item(198, 4).Value = 46
$ws.Cells.Item(198, 5).Value = 0
$ws.Cells.Item(198, 8).Value = 0

# Row 214 (now Montserrat)
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

# Row 215 (now Islas Malvinas)
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
